$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reviews in rows 7 and 8 ("confirm" column G) flip from "no" to "yes"
$ws.Range("G7").Value = "yes"
$ws.Range("G8").Value = "yes"

# Move the active selection/cursor to G9 (was G10)
$ws.Range("G9").Select()
